$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Prix Spot" (sheet1): add new column AN (23-jul) with header
# in row 1 and 24 hourly values in rows 2-25.
# -----------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (AM1) onto the new
# header cell (AN1) so it keeps the bold/centered/bordered style, then
# set its text.
$wsSpot.Range("AM1").Copy()
$wsSpot.Range("AN1").PasteSpecial(-4122)
$wsSpot.Range("AN1").Value = "23-jul"

$spotValues = @{
    2  = 89.79000000000001
    3  = 82.64
    4  = 80.94
    5  = 78.26000000000001
    6  = 77.83
    7  = 81.11
    8  = 90.87
    9  = 100.5
    10 = 96.5
    11 = 83.31
    12 = 76.40000000000001
    13 = 68.64
    14 = 61.97
    15 = 51.89
    16 = 44.88
    17 = 50.38
    18 = 71.63
    19 = 82.17
    20 = 92.17
    21 = 104.98
    22 = 111.63
    23 = 112
    24 = 109.96
    25 = 100.39
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 40).Value = $spotValues[$row]
}

# -----------------------------------------------------------------
# Sheet "Gaz" (sheet2): append a new row 37 with date 2025-07-21 and
# its price.
# -----------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A37").NumberFormat = "@"
$wsGaz.Range("A37").Value = "2025-07-21"
$wsGaz.Range("A37").Style = "Normal"
$wsGaz.Range("B37").Value = 32.6

# -----------------------------------------------------------------
# Sheet "CO2" (sheet3): append a new row 37 with date 2025-07-21 and
# its price.
# -----------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A37").NumberFormat = "@"
$wsCo2.Range("A37").Value = "2025-07-21"
$wsCo2.Range("A37").Style = "Normal"
$wsCo2.Range("B37").Value = 69.09999999999999
